# Update "想去人数" (F column) figures for multiple sheets, reflecting
# output regenerated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 15072
$ws1.Range("F3").Value = 19216
$ws1.Range("F5").Value = 146
$ws1.Range("F6").Value = 56
$ws1.Range("F14").Value = 173
$ws1.Range("F22").Value = 8026
$ws1.Range("F27").Value = 1252
$ws1.Range("F29").Value = 6084
$ws1.Range("F30").Value = 118
$ws1.Range("F31").Value = 74
$ws1.Range("F35").Value = 5486
$ws1.Range("F36").Value = 996

# --- Sheet: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 19

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 15072
$ws4.Range("F3").Value = 19216
$ws4.Range("F5").Value = 146
$ws4.Range("F6").Value = 56
$ws4.Range("F14").Value = 173
$ws4.Range("F23").Value = 8026
$ws4.Range("F28").Value = 1252
$ws4.Range("F30").Value = 19
$ws4.Range("F32").Value = 6084
$ws4.Range("F33").Value = 118
$ws4.Range("F34").Value = 74
$ws4.Range("F38").Value = 5486
$ws4.Range("F39").Value = 996
